$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (was "Hoja1")
$ws.Name = "visorInformacionTecnicaRed"

# Move/change the active selection from G6 to D5
$ws.Range("D5").Select() | Out-Null

# Give the header row (row 1) an explicit height of 30
$ws.Rows.Item(1).RowHeight = 30
